# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly refreshed figures, as published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rows are shifted by -1 relative to "全部类型")
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value  = 23
$wsExhibit.Range("F8").Value  = 13841
$wsExhibit.Range("F10").Value = 76
$wsExhibit.Range("F11").Value = 5587
$wsExhibit.Range("F13").Value = 41
$wsExhibit.Range("F19").Value = 749
$wsExhibit.Range("F21").Value = 38
$wsExhibit.Range("F22").Value = 10379
$wsExhibit.Range("F24").Value = 22
$wsExhibit.Range("F25").Value = 37
$wsExhibit.Range("F26").Value = 3699

# Sheet "全部类型" (same events, one row lower)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 23
$wsAll.Range("F9").Value  = 13841
$wsAll.Range("F11").Value = 76
$wsAll.Range("F12").Value = 5587
$wsAll.Range("F14").Value = 41
$wsAll.Range("F20").Value = 749
$wsAll.Range("F22").Value = 38
$wsAll.Range("F24").Value = 10379
$wsAll.Range("F26").Value = 22
$wsAll.Range("F27").Value = 37
$wsAll.Range("F28").Value = 3699
